$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Paragraph 1: "Yerba Buena, 17 de Diciembre de 2013"
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Format.KeepWithNext = $true
$p1.Format.SpaceAfter = 12

# ---------------------------------------------------------------------
# Paragraph 2: "ORDENANZA Nº 1949"
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$p2.Format.KeepWithNext = $true
$p2.Format.SpaceBefore = 12
$p2.Format.SpaceAfter = 18
$p2.Range.Font.Bold = 1

# ---------------------------------------------------------------------
# Paragraph 3: "EL CONCEJO DELIBERANTE SANCIONA CON FUERZA DE ORDENANZA"
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3.Format.KeepWithNext = $true
$p3.Format.SpaceBefore = 18
$p3.Format.SpaceAfter = 18
$p3.Format.LeftIndent = 99.2
$p3.Format.RightIndent = 99.2
$p3.Range.Font.Bold = 1

# ---------------------------------------------------------------------
# Paragraph 4: "ARTICULO PRIMERO: REFRENDASE el Decreto Nº 496, ..."
# ---------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$p4.Format.KeepWithNext = $true
$p4.Format.SpaceAfter = 6
$p4.Format.Alignment = 0

# Underline "ARTICULO PRIMERO:" (word + colon), then drop the leading
# space off the run that used to read ": " so it is its own run.
$rng = $d.Content
$rng.Find.Execute("ARTICULO PRIMERO:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$wordRange = $d.Range($rng.Start, $rng.End - 1)
$wordRange.Font.Underline = 1
$colonRange = $d.Range($rng.End - 1, $rng.End)
$colonRange.Font.Underline = 1

# Split "ecreto Nº 496" -> "ecreto N" + "º " + "496" (three plain runs).
$rng2 = $d.Content
$rng2.Find.Execute("ecreto N", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitRange = $d.Range($rng2.End, $rng2.End + 2)
$splitRange.Font.Bold = 1
$splitRange.Font.Bold = 0

# ---------------------------------------------------------------------
# Paragraph 5: "ARTICULO SEGUNDO: COMUNIQUESE, REGISTRESE Y ARCHIVESE."
# ---------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$p5.Format.KeepWithNext = $true
$p5.Format.SpaceAfter = 6
$p5.Format.Alignment = 0

# Remove the stray leading space before "ARTICULO SEGUNDO".
$rng3 = $d.Content
$rng3.Find.Execute("ARTICULO SEGUNDO:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$leadSpace = $d.Range($rng3.Start - 1, $rng3.Start)
$leadSpace.Delete()

# Underline "ARTICULO SEGUNDO:" (word + colon).
$rng4 = $d.Content
$rng4.Find.Execute("ARTICULO SEGUNDO:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$wordRange2 = $d.Range($rng4.Start, $rng4.End - 1)
$wordRange2.Font.Underline = 1
$colonRange2 = $d.Range($rng4.End - 1, $rng4.End)
$colonRange2.Font.Underline = 1

# ---------------------------------------------------------------------
# Section properties: starting page number
# ---------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$sec.Headers.Item(1).PageNumbers.StartingNumber = 2796

# ---------------------------------------------------------------------
# Default footer (footer2.xml): collapse the three paragraphs into one
# and restyle the surviving run properties.
# ---------------------------------------------------------------------
$ftr = $sec.Footers.Item(1)
$ftr.Range.Paragraphs.Item(2).Range.Delete()
$ftr.Range.Paragraphs.Item(1).Range.Delete()
$ftrPara = $ftr.Range.Paragraphs.Item(1)
$ftrPara.Range.Font.NameAscii = "Book Antiqua"
$ftrPara.Range.Font.NameOther = "Book Antiqua"
$ftrPara.Range.Font.Size = 10
$ftrPara.Range.Font.Color = 8421504
